$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-18 22:48:52'
$ws.Range('I2').Value = '1.4 mm'
$ws.Range('O2').Value = '1.8 °C'
$ws.Range('E3').Value = '2026-02-18 22:48:54'
$ws.Range('E4').Value = '2026-02-18 22:48:57'
$ws.Range('J4').Value = '1012.3 hPa'
$ws.Range('E5').Value = '2026-02-18 22:49:00'
$ws.Range('I5').Value = '1.4 mm'
$ws.Range('N5').Value = '-3.8 °C 22:29 TU'
$ws.Range('O5').Value = '0.5 °C'
$ws.Range('E6').Value = '2026-02-18 22:49:03'
$ws.Range('J6').Value = '1012.0 hPa'
$ws.Range('E7').Value = '2026-02-18 22:49:05'
$ws.Range('J7').Value = '1013.6 hPa'
$ws.Range('E8').Value = '2026-02-18 22:49:07'
$ws.Range('J8').Value = '1013.3 hPa'
$ws.Range('E9').Value = '2026-02-18 22:49:10'
$ws.Range('E10').Value = '2026-02-18 22:49:13'
$ws.Range('E11').Value = '2026-02-18 22:49:16'
$ws.Range('E12').Value = '2026-02-18 22:49:19'
$ws.Range('O12').Value = '11.1 °C'
$ws.Range('E13').Value = '2026-02-18 22:49:21'
$ws.Range('J13').Value = '1014.6 hPa'
$ws.Range('O13').Value = '4.1 °C'
$ws.Range('E14').Value = '2026-02-18 22:49:23'
$ws.Range('E15').Value = '2026-02-18 22:49:26'
$ws.Range('E16').Value = '2026-02-18 22:49:28'
$ws.Range('G16').Value = '74 cm'
$ws.Range('I16').Value = '2.7 mm'
$ws.Range('N16').Value = '-3.5 °C 22:14 TU'
$ws.Range('O16').Value = '-0.3 °C'
$ws.Range('E17').Value = '2026-02-18 22:49:31'
$ws.Range('E18').Value = '2026-02-18 22:49:34'
$ws.Range('J18').Value = '1012.5 hPa'
$ws.Range('E19').Value = '2026-02-18 22:49:37'
$ws.Range('E20').Value = '2026-02-18 22:49:39'
$ws.Range('I20').Value = '1.0 mm'
$ws.Range('N20').Value = '-3.4 °C 22:29 TU'
$ws.Range('O20').Value = '-0.6 °C'
$ws.Range('E21').Value = '2026-02-18 22:49:42'
$ws.Range('J21').Value = '1014.1 hPa'
$ws.Range('O21').Value = '6.7 °C'
$ws.Range('E22').Value = '2026-02-18 22:49:45'
$ws.Range('N22').Value = '-6.2 °C 22:17 TU'
$ws.Range('O22').Value = '-1.9 °C'
$ws.Range('E23').Value = '2026-02-18 22:49:48'
$ws.Range('I23').Value = '0.8 mm'
$ws.Range('O23').Value = '-0.2 °C'
$ws.Range('E24').Value = '2026-02-18 22:49:50'
$ws.Range('J24').Value = '1014.2 hPa'
$ws.Range('E25').Value = '2026-02-18 22:49:53'
$ws.Range('I25').Value = '0.1 mm'
$ws.Range('O25').Value = '1.6 °C'
$ws.Range('E26').Value = '2026-02-18 22:49:56'
$ws.Range('J26').Value = '1011.5 hPa'
$ws.Range('E27').Value = '2026-02-18 22:49:59'
$ws.Range('I27').Value = '0.6 mm'
$ws.Range('N27').Value = '-1.7 °C 22:28 TU'
$ws.Range('O27').Value = '1.3 °C'
$ws.Range('E28').Value = '2026-02-18 22:50:02'
$ws.Range('J28').Value = '1012.2 hPa'
$ws.Range('E29').Value = '2026-02-18 22:50:04'
$ws.Range('E30').Value = '2026-02-18 22:50:07'
$ws.Range('J30').Value = '1011.7 hPa'
$ws.Range('E31').Value = '2026-02-18 22:50:10'
$ws.Range('J31').Value = '1010.6 hPa'
$ws.Range('O31').Value = '12.5 °C'
$ws.Range('E32').Value = '2026-02-18 22:50:13'
$ws.Range('L32').Value = '34.9 km/h - 239º 22:15 TU'
$ws.Range('E33').Value = '2026-02-18 22:50:15'
$ws.Range('J33').Value = '1013.4 hPa'
$ws.Range('E34').Value = '2026-02-18 22:50:18'
$ws.Range('I34').Value = '0.2 mm'
$ws.Range('E35').Value = '2026-02-18 22:50:21'
$ws.Range('I35').Value = '0.1 mm'
$ws.Range('J35').Value = '1013.9 hPa'
$ws.Range('L35').Value = '46.1 km/h - 278º 22:04 TU'
$ws.Range('E36').Value = '2026-02-18 22:50:24'
$ws.Range('J36').Value = '1012.2 hPa'
$ws.Range('L36').Value = '46.8 km/h - 1º 22:29 TU'
$ws.Range('E37').Value = '2026-02-18 22:50:26'
$ws.Range('J37').Value = '1013.8 hPa'
$ws.Range('E38').Value = '2026-02-18 22:50:29'
$ws.Range('E39').Value = '2026-02-18 22:50:32'
$ws.Range('O39').Value = '1.1 °C'
$ws.Range('E40').Value = '2026-02-18 22:50:34'
$ws.Range('J40').Value = '1014.7 hPa'
$ws.Range('E41').Value = '2026-02-18 22:50:37'
$ws.Range('J41').Value = '1013.9 hPa'
$ws.Range('E42').Value = '2026-02-18 22:50:40'
$ws.Range('E43').Value = '2026-02-18 22:50:42'
$ws.Range('E44').Value = '2026-02-18 22:50:45'
$ws.Range('O44').Value = '-1.6 °C'
$ws.Range('E45').Value = '2026-02-18 22:50:48'
$ws.Range('J45').Value = '1011.4 hPa'
$ws.Range('O45').Value = '7.2 °C'
$ws.Range('E46').Value = '2026-02-18 22:50:51'
$ws.Range('J46').Value = '1014.3 hPa'

# Percentage-looking strings need special handling to avoid Excel auto-converting
# them to numeric percentages; force text format, then restore original style via copy/paste.
$ws.Range('H15').NumberFormat = '@'
$ws.Range('H15').Value = '79%'
$ws.Range('C2').Copy() | Out-Null
$ws.Range('H15').PasteSpecial(-4122) | Out-Null
$ws.Range('H18').NumberFormat = '@'
$ws.Range('H18').Value = '78%'
$ws.Range('C2').Copy() | Out-Null
$ws.Range('H18').PasteSpecial(-4122) | Out-Null
$ws.Range('H19').NumberFormat = '@'
$ws.Range('H19').Value = '86%'
$ws.Range('C2').Copy() | Out-Null
$ws.Range('H19').PasteSpecial(-4122) | Out-Null
$ws.Range('H20').NumberFormat = '@'
$ws.Range('H20').Value = '77%'
$ws.Range('C2').Copy() | Out-Null
$ws.Range('H20').PasteSpecial(-4122) | Out-Null
$ws.Range('H23').NumberFormat = '@'
$ws.Range('H23').Value = '60%'
$ws.Range('C2').Copy() | Out-Null
$ws.Range('H23').PasteSpecial(-4122) | Out-Null
$ws.Range('H25').NumberFormat = '@'
$ws.Range('H25').Value = '51%'
$ws.Range('C2').Copy() | Out-Null
$ws.Range('H25').PasteSpecial(-4122) | Out-Null
$ws.Range('H27').NumberFormat = '@'
$ws.Range('H27').Value = '59%'
$ws.Range('C2').Copy() | Out-Null
$ws.Range('H27').PasteSpecial(-4122) | Out-Null
$ws.Range('H29').NumberFormat = '@'
$ws.Range('H29').Value = '85%'
$ws.Range('C2').Copy() | Out-Null
$ws.Range('H29').PasteSpecial(-4122) | Out-Null
$ws.Range('H34').NumberFormat = '@'
$ws.Range('H34').Value = '49%'
$ws.Range('C2').Copy() | Out-Null
$ws.Range('H34').PasteSpecial(-4122) | Out-Null
$ws.Range('H39').NumberFormat = '@'
$ws.Range('H39').Value = '43%'
$ws.Range('C2').Copy() | Out-Null
$ws.Range('H39').PasteSpecial(-4122) | Out-Null
$ws.Range('H45').NumberFormat = '@'
$ws.Range('H45').Value = '66%'
$ws.Range('C2').Copy() | Out-Null
$ws.Range('H45').PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
